$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the payment-date columns to use a live TODAY() formula instead of
# a hard-coded date serial, then leave the selection on the edited cell.
$ws.Range("D2").Formula = "=TODAY()"
$ws.Range("E2").Formula = "=TODAY()"

$ws.Range("D2").Select() | Out-Null
